$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.309.05'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").Value = '1.868.07'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("D4").Value = '''0.9997'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''235.40'
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("D6").Value = '''0.9996'
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").Value = '''0.2842'
$ws.Range("E8").Value = '  +0.37%  '
$ws.Range("D9").Value = '''0.06533'
$ws.Range("E9").Value = '  -0.87%  '
$ws.Range("D10").Value = '''21.40'
$ws.Range("E10").Value = '  +3.87%  '
$ws.Range("D11").Value = '''0.07874'
$ws.Range("E11").Value = '  +1.39%  '
$ws.Range("D12").Value = '''97.86'
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").Value = '1.873.34'
$ws.Range("E13").Value = '  -0.38%  '
$ws.Range("D14").Value = '''5.104'
$ws.Range("E14").Value = '  +0.63%  '
$ws.Range("D15").Value = '''0.6755'
$ws.Range("E15").Value = '  +0.21%  '
$ws.Range("D16").Value = '''279.64'
$ws.Range("E16").Value = '  -1.68%  '
$ws.Range("D17").Value = '30.293.13'
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").Value = '''0.9999'
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("D19").Value = '''5.500'
$ws.Range("E19").Value = '  +1.86%  '
$ws.Range("D20").Value = '''12.71'
$ws.Range("E20").Value = '  +0.76%  '
$ws.Range("D21").Value = '2.111.09'
$ws.Range("E21").Value = '  -0.90%  '
$ws.Range("D22").Value = '''0.000007287'
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").Value = '''0.9995'
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("D24").Value = '''6.165'
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").Value = '''9.193'
$ws.Range("E25").Value = '  -1.92%  '
$ws.Range("D26").Value = '''165.12'
$ws.Range("E26").Value = '  -1.84%  '
$ws.Range("D27").Value = '''19.13'
$ws.Range("E27").Value = '  -0.28%  '
$ws.Range("D28").Value = '''1.929'
$ws.Range("E28").Value = '  -2.93%  '
$ws.Range("D29").Value = '''1.377'
$ws.Range("E29").Value = '  +0.43%  '
$ws.Range("E30").Value = '  -0.27%  '
$ws.Range("D31").Value = '''4.378'
$ws.Range("E31").Value = '  +0.14%  '
$ws.Range("D32").Value = '''1.477'
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("D33").Value = '''4.100'
$ws.Range("E33").Value = '  -0.54%  '
$ws.Range("D34").Value = '''0.04714'
$ws.Range("E34").Value = '  +0.94%  '
$ws.Range("D35").Value = '''1.131'
$ws.Range("E35").Value = '  +3.10%  '
$ws.Range("D36").Value = '''0.7063'
$ws.Range("D37").Value = '''2.726'
$ws.Range("E37").Value = '  +0.49%  '
$ws.Range("D38").Value = '''0.01858'
$ws.Range("E38").Value = '  -0.58%  '
$ws.Range("D39").Value = '''6.254'
$ws.Range("E39").Value = '  -5.07%  '
$ws.Range("D40").Value = '''2.529'
$ws.Range("E40").Value = '  +0.26%  '
$ws.Range("D41").Value = '''73.50'
$ws.Range("E42").Value = '  -0.89%  '
$ws.Range("D43").Value = '''0.8483'
$ws.Range("E43").Value = '  -1.96%  '
$ws.Range("D44").Value = '''0.4178'
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("D45").Value = '''0.9996'
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").Value = '''103.70'
$ws.Range("E46").Value = '  +0.74%  '
$ws.Range("D47").Value = '''7.179'
$ws.Range("E47").Value = '  -1.20%  '
$ws.Range("D48").Value = '''9.248'
$ws.Range("E48").Value = '  +0.68%  '
$ws.Range("D49").Value = '''936.29'
$ws.Range("E49").Value = '  -4.69%  '
$ws.Range("D50").Value = '''34.08'
$ws.Range("E50").Value = '  +0.73%  '
$ws.Range("D51").Value = '''0.1124'
$ws.Range("E51").Value = '  -1.85%  '
